$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (student 48951): Penalty recalculated
$ws.Range("F2").Value = 86.49093270617772

# Row 3 (student 48953): crossed into "Good" tier
$ws.Range("B3").Value = 80
$ws.Range("C3").Value = 60.64776531782864
$ws.Range("D3").Value = -19.35223468217136
$ws.Range("F3").Value = 19.35223468217136
$ws.Range("G3").Value = "Good"

# Row 4 (student 48955): Penalty recalculated
$ws.Range("F4").Value = 86.49093270617772

# Row 5 (student 48956)
$ws.Range("B5").Value = 78
$ws.Range("C5").Value = 67.02237126211561
$ws.Range("D5").Value = -10.97762873788439
$ws.Range("F5").Value = 10.97762873788439

# Row 6 (student 48957)
$ws.Range("B6").Value = 78
$ws.Range("C6").Value = 67.02237126211561
$ws.Range("D6").Value = -10.97762873788439
$ws.Range("F6").Value = 10.97762873788439

# Row 7 (student 48961): Excellence tier
$ws.Range("B7").Value = 91.5
$ws.Range("C7").Value = 84.66243092374654
$ws.Range("D7").Value = -6.837569076253459
$ws.Range("F7").Value = 6.837569076253453
$ws.Range("G7").Value = "Excellence"

# Row 8 (student 48962)
$ws.Range("B8").Value = 66.5
$ws.Range("C8").Value = 32.96269956291068
$ws.Range("D8").Value = -33.53730043708932
$ws.Range("F8").Value = 33.53730043708932

# Row 10 (student 48970)
$ws.Range("C10").Value = 55.3506518718915
$ws.Range("D10").Value = -25.6493481281085
$ws.Range("F10").Value = 25.6493481281085

# Row 11 (student 48971)
$ws.Range("B11").Value = 53
$ws.Range("C11").Value = 15.04705546591035
$ws.Range("D11").Value = -37.95294453408965
$ws.Range("F11").Value = 37.95294453408965

# Row 12 (student 48972)
$ws.Range("B12").Value = 71
$ws.Range("C12").Value = 31.85099496236071
$ws.Range("D12").Value = -39.14900503763929
$ws.Range("F12").Value = 39.14900503763929

# Row 13 (student 48973): Penalty recalculated
$ws.Range("F13").Value = 45.14824869927797

# Row 14 (student 48975)
$ws.Range("C14").Value = 55.3506518718915
$ws.Range("D14").Value = -25.6493481281085
$ws.Range("F14").Value = 25.6493481281085

# Row 15 (student 48976)
$ws.Range("C15").Value = 21.85298835372897
$ws.Range("D15").Value = -41.64701164627103
$ws.Range("F15").Value = 41.64701164627103

# Row 16 (student 48978)
$ws.Range("B16").Value = 67.5
$ws.Range("C16").Value = 23.62611504402493
$ws.Range("D16").Value = -43.87388495597507
$ws.Range("F16").Value = 43.87388495597507

# Row 17 (student 48980)
$ws.Range("B17").Value = 53
$ws.Range("D17").Value = -53
$ws.Range("F17").Value = 63.44838747479471

# Row 19 (student 48987): Penalty recalculated
$ws.Range("F19").Value = 45.14824869927797

# Row 20 (student 48991)
$ws.Range("B20").Value = 67.5
$ws.Range("C20").Value = 23.62611504402493
$ws.Range("D20").Value = -43.87388495597507
$ws.Range("F20").Value = 43.87388495597507

# Row 21 (student 48993)
$ws.Range("B21").Value = 71
$ws.Range("C21").Value = 31.85099496236071
$ws.Range("D21").Value = -39.14900503763929
$ws.Range("F21").Value = 39.14900503763929

# Row 22 (student 48994)
$ws.Range("B22").Value = 71
$ws.Range("C22").Value = 31.85099496236071
$ws.Range("D22").Value = -39.14900503763929
$ws.Range("F22").Value = 39.14900503763929

# Row 23 (student 74309)
$ws.Range("C23").Value = 21.85298835372897
$ws.Range("D23").Value = -41.64701164627103
$ws.Range("F23").Value = 41.64701164627103

# Row 24 (student 74311): moved from Excellence(92.5) to Excellence(91.5), adjusted
$ws.Range("B24").Value = 91.5
$ws.Range("C24").Value = 84.66243092374654
$ws.Range("D24").Value = -6.837569076253459
$ws.Range("F24").Value = 6.837569076253453

# Materialize an empty row 25 marker (touch a row-level property with its
# existing default value so the row element is written without attributes
# or cells), matching the blank <row r="25"/> introduced by the edit.
$ws.Rows.Item(25).OutlineLevel = 0
